$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.948.08'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.50%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.027.68'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.32%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.06'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '153.07'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.26%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.024.52'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.33%  '
$ws.Range("E9").Value = '  -0.22%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.06'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +16.44%  '
$ws.Range("E11").Value = '  +1.89%  '
$ws.Range("E12").Value = '  +2.10%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000233'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.93%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.63'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.37%  '
$ws.Range("E15").Value = '  -1.41%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.529.69'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.37%  '
$ws.Range("E17").Value = '  +2.59%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '62.958.05'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.028.09'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.33%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '448.53'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.09%  '
$ws.Range("E21").Value = '  +1.60%  '
$ws.Range("E22").Value = '  +1.91%  '
$ws.Range("E23").Value = '  +2.65%  '
$ws.Range("E24").Value = '  +7.84%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.09'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.62%  '
$ws.Range("E26").Value = '  +6.55%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.38'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.13%  '
$ws.Range("E28").Value = '  +0.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.54'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.68%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.28'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +11.11%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.69'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.52%  '
$ws.Range("E32").Value = '  +0.00%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.64'
$ws.Range("D33").Style = "Normal"
$ws.Range("E34").Value = '  +1.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0₃0875'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.84%  '
$ws.Range("E36").Value = '  +3.09%  '
$ws.Range("E37").Value = '  +1.73%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.16'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +10.29%  '
$ws.Range("E39").Value = '  +2.96%  '
$ws.Range("E40").Value = '  +8.43%  '
$ws.Range("E41").Value = '  +0.68%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.05'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '44.60'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +17.01%  '
$ws.Range("E44").Value = '  +14.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '392.60'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.56%  '
$ws.Range("E46").Value = '  +2.17%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.709.52'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.60%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '133.61'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.43%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '26.47'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +13.40%  '
$ws.Range("E51").Value = '  +5.87%  '
